{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items,text\");\nawait context.sync();\n\nfor (let i = paras.items.length - 1; i >= 0; i--) {\n  const t = paras.items[i].text;\n  if (t.indexOf(\"Diagrammi delle Attivit\u00e0\") !== -1) {\n    paras.items[i].delete();\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$found = $false\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Diagrammi delle Attivit\u00e0*\") {\n        $p.Range.Delete()\n        $found = $true\n        break\n    }\n}\n"}
